$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column O ("LD Target Qty/Day") values per processed executive data (test2.py)
$ws.Cells.Item(2, 15).Value = 160
$ws.Cells.Item(3, 15).Value = 9
$ws.Cells.Item(4, 15).Value = 165
$ws.Cells.Item(6, 15).Value = 9
$ws.Cells.Item(7, 15).Value = 405
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(13, 15).Value = 496
$ws.Cells.Item(15, 15).Value = 106
$ws.Cells.Item(16, 15).Value = 811
$ws.Cells.Item(17, 15).Value = 541
$ws.Cells.Item(18, 15).Value = 899
$ws.Cells.Item(19, 15).Value = 53
$ws.Cells.Item(20, 15).Value = 20
$ws.Cells.Item(21, 15).Value = 787
$ws.Cells.Item(22, 15).Value = 82
$ws.Cells.Item(23, 15).Value = 653
$ws.Cells.Item(24, 15).Value = 118
$ws.Cells.Item(25, 15).Value = 552
$ws.Cells.Item(26, 15).Value = 345
$ws.Cells.Item(27, 15).Value = 5
$ws.Cells.Item(28, 15).Value = 1
$ws.Cells.Item(29, 15).Value = 30
$ws.Cells.Item(30, 15).Value = 6
$ws.Cells.Item(31, 15).Value = 13
$ws.Cells.Item(32, 15).Value = 6
$ws.Cells.Item(33, 15).Value = 14
$ws.Cells.Item(34, 15).Value = 0
$ws.Cells.Item(35, 15).Value = 283
$ws.Cells.Item(36, 15).Value = 87
$ws.Cells.Item(37, 15).Value = 115
$ws.Cells.Item(38, 15).Value = 163
$ws.Cells.Item(39, 15).Value = 56
$ws.Cells.Item(40, 15).Value = 338
$ws.Cells.Item(42, 15).Value = 31
$ws.Cells.Item(43, 15).Value = 4
$ws.Cells.Item(44, 15).Value = 198
$ws.Cells.Item(46, 15).Value = 624
$ws.Cells.Item(47, 15).Value = 240
$ws.Cells.Item(48, 15).Value = 221
$ws.Cells.Item(49, 15).Value = 44
$ws.Cells.Item(50, 15).Value = 7
$ws.Cells.Item(51, 15).Value = 7
$ws.Cells.Item(52, 15).Value = 9
$ws.Cells.Item(53, 15).Value = 223
$ws.Cells.Item(54, 15).Value = 6
$ws.Cells.Item(55, 15).Value = 7
$ws.Cells.Item(56, 15).Value = 214
$ws.Cells.Item(57, 15).Value = 79
$ws.Cells.Item(58, 15).Value = 377
$ws.Cells.Item(60, 15).Value = 0
$ws.Cells.Item(61, 15).Value = 24
$ws.Cells.Item(62, 15).Value = 355
$ws.Cells.Item(63, 15).Value = 30
$ws.Cells.Item(64, 15).Value = 110
$ws.Cells.Item(65, 15).Value = 6
$ws.Cells.Item(66, 15).Value = 88
$ws.Cells.Item(67, 15).Value = 54
$ws.Cells.Item(68, 15).Value = 64
$ws.Cells.Item(70, 15).Value = 142
$ws.Cells.Item(72, 15).Value = 2000
$ws.Cells.Item(73, 15).Value = 30
$ws.Cells.Item(74, 15).Value = 104
$ws.Cells.Item(77, 15).Value = 143
$ws.Cells.Item(78, 15).Value = 685
$ws.Cells.Item(79, 15).Value = 125
$ws.Cells.Item(80, 15).Value = 759
$ws.Cells.Item(81, 15).Value = 674
